# Applies the diff: rewrites rows 2-8 with the new deduplicated log data
# and removes the old trailing rows 9-10 (dimension A1:H10 -> A1:H8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-absent rows (old rows 9 and 10).
$ws.Rows.Item(9).EntireRow.Delete()
$ws.Rows.Item(9).EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "life-dev/main"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "scroll"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "channel, page_url, scroll_rate, os_name"
$ws.Range("G2").Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 75, iOS"
$ws.Range("H2").Value = "4"

# Row 3
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "life-dev/main"
$ws.Range("C3").Value = "상품"
$ws.Range("D3").Value = "click"
$ws.Range("E3").Value = "(스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P"
$ws.Range("F3").Value = "channel, page_url, click_text, module_id, module_order, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_is_ad, el_order, module_name, os_name"
$ws.Range("G3").Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, C-3, 13, 1, 3086, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, 마이그스토어, 40,000원, 20,000원, 50%, F, 1, commerce-category-ranking, iOS"
$ws.Range("H3").Value = "16"

# Row 4
$ws.Range("A4").Value = "5"
$ws.Range("B4").Value = "life-dev/main"
$ws.Range("C4").Value = "뉴스"
$ws.Range("D4").Value = "click"
$ws.Range("E4").Value = "[OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화"
$ws.Range("F4").Value = "channel, page_url, click_text, module_id, module_order, el_order, module_name, article_title, os_name"
$ws.Range("G4").Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, [OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화, D-1, 14, 1, news-card, [OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화, iOS"
$ws.Range("H4").Value = "9"

# Row 5
$ws.Range("A5").Value = "6"
$ws.Range("B5").Value = "life-dev/news/detail/10736"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "click"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "channel, page_url, os_name"
$ws.Range("G5").Value = "Rround, https://life-dev.hectoinnovation.co.kr/news/detail/10736, iOS"
$ws.Range("H5").Value = "3"

# Row 6
$ws.Range("A6").Value = "11"
$ws.Range("B6").Value = "ecommerce-dev/product/detail/3086"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "pageview"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "channel, page_url, prd_code, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name"
$ws.Range("G6").Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/3086, 3086, 40,000원, 18,000원, 55%, 0, 0, #포도씨유___#올리브유___#소르바스___#압착오일___#엑스트라버진___#해바라기유___#카놀라유___#유기농, iOS"
$ws.Range("H6").Value = "10"

# Row 7
$ws.Range("A7").Value = "12"
$ws.Range("B7").Value = "ecommerce-dev/product/detail/3086"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "click"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "channel, page_url, tab_name, prd_code, prd_name, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name"
$ws.Range("G7").Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/3086, 상품상세`n, 3086, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, 40,000원, 20,000원, 55%, 0, 0, #포도씨유___#올리브유___#소르바스___#압착오일___#엑스트라버진___#해바라기유___#카놀라유___#유기농, iOS"
$ws.Range("H7").Value = "12"

# Row 8
$ws.Range("A8").Value = "14"
$ws.Range("B8").Value = "life-dev/main"
$ws.Range("C8").Value = "뉴스"
$ws.Range("D8").Value = "click"
$ws.Range("E8").Value = "K뮤지컬 통했다…'어쩌면 해피엔딩', 토니상 극본상·음악상 수상"
$ws.Range("F8").Value = "channel, page_url, click_text, module_id, module_order, el_order, module_name, article_title, os_name"
$ws.Range("G8").Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, K뮤지컬 통했다…'어쩌면 해피엔딩', 토니상 극본상·음악상 수상, D-1, 19, 1, news-card, K뮤지컬 통했다…'어쩌면 해피엔딩', 토니상 극본상·음악상 수상, iOS"
$ws.Range("H8").Value = "9"
